$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark that currently sits in the title
#    paragraph (between "MP73010" and " - Assignment 1 exercise").
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Locate the first of the two trailing empty paragraphs (right after the
#    "Ben changing things up!" paragraph) and replace it with:
#      - a new paragraph containing the "Hello,it's  great using github" text
#        (with the spell/grammar-checker proofErr markers Word leaves behind)
#      - a new paragraph that now carries the "_GoBack" bookmark
#    The second trailing empty paragraph is left untouched.
# ---------------------------------------------------------------------------
$apos = [char]0x2019

$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -eq [char]13 -and $cand.Range.Text.Length -eq 1) {
        $targetPara = $cand
        break
    }
}

$newBodyXml = `
    '<w:p>' + `
        '<w:proofErr w:type="spellStart"/>' + `
        '<w:r><w:t>Hello</w:t></w:r>' + `
        '<w:proofErr w:type="gramStart"/>' + `
        '<w:r><w:t>,it' + $apos + 's</w:t></w:r>' + `
        '<w:proofErr w:type="spellEnd"/>' + `
        '<w:proofErr w:type="gramEnd"/>' + `
        '<w:r><w:t xml:space="preserve">  great using </w:t></w:r>' + `
        '<w:proofErr w:type="spellStart"/>' + `
        '<w:r><w:t>github</w:t></w:r>' + `
        '<w:proofErr w:type="spellEnd"/>' + `
    '</w:p>' + `
    '<w:p>' + `
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
        '<w:bookmarkEnd w:id="0"/>' + `
    '</w:p>'

$pkgXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
            '<pkg:xmlData>' + `
                '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
                    '<w:body>' + $newBodyXml + '</w:body>' + `
                '</w:document>' + `
            '</pkg:xmlData>' + `
        '</pkg:part>' + `
    '</pkg:package>'

$targetPara.Range.InsertXML($pkgXml)

Write-Output "Done. Paragraph count now: $($d.Paragraphs.Count)"
